# Vega Monumental Concepción - Zapallo: weekly update
# Insert two new weekly price rows at the top of the "Camote" block (row 224)
# which shifts the existing rows 224:255 down to 226:257.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 224 (inserted twice shifts
# everything else down by 2 rows total).
$ws.Rows.Item(224).Insert()
$ws.Rows.Item(224).Insert()

# New row 224: Camote 1a (guarda), 2022-08-25
$ws.Cells.Item(224, 1).Value = 11
$ws.Cells.Item(224, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(224, 3).Value = "Bíobío"
$ws.Cells.Item(224, 4).Value = 44798
$ws.Cells.Item(224, 5).Value = 8
$ws.Cells.Item(224, 6).Value = 100112045
$ws.Cells.Item(224, 7).Value = "Zapallo"
$ws.Cells.Item(224, 8).Value = "Camote"
$ws.Cells.Item(224, 9).Value = "1a (guarda)"
$ws.Cells.Item(224, 10).Value = 600
$ws.Cells.Item(224, 11).Value = 750
$ws.Cells.Item(224, 12).Value = 800
$ws.Cells.Item(224, 13).Value = 775
$ws.Cells.Item(224, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(224, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(224, 16).Value = 775
$ws.Cells.Item(224, 17).Value = 1
$ws.Cells.Item(224, 18).Value = "Hortaliza"

# New row 225: Camote 2a (guarda), 2022-08-25
$ws.Cells.Item(225, 1).Value = 11
$ws.Cells.Item(225, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(225, 3).Value = "Bíobío"
$ws.Cells.Item(225, 4).Value = 44798
$ws.Cells.Item(225, 5).Value = 8
$ws.Cells.Item(225, 6).Value = 100112045
$ws.Cells.Item(225, 7).Value = "Zapallo"
$ws.Cells.Item(225, 8).Value = "Camote"
$ws.Cells.Item(225, 9).Value = "2a (guarda)"
$ws.Cells.Item(225, 10).Value = 300
$ws.Cells.Item(225, 11).Value = 650
$ws.Cells.Item(225, 12).Value = 650
$ws.Cells.Item(225, 13).Value = 650
$ws.Cells.Item(225, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(225, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(225, 16).Value = 650
$ws.Cells.Item(225, 17).Value = 1
$ws.Cells.Item(225, 18).Value = "Hortaliza"

Write-Host "done"
